$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").Value = "ValueCH"
$ws.Range("F2").Select()
